$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fractions")

# Fix header row labels: the columns were mislabeled.
$ws.Range("A1").Value = "TMT_Set"
$ws.Range("B1").Value = "LCMS_Injection"
$ws.Range("C1").Value = "Fraction"

# Move selection to C5 (as recorded in the saved view state).
$ws.Range("C5").Select()
